# Drop in all data files from 3.0 RMI script
#
# This workbook previously contained a scratch "Texas Notes" sheet that was
# used while the BAU Fraction of Components Sold by Quality Level data was
# being adapted for Texas. This change removes that scratch sheet and
# restores the underlying "Data" sheet values for Gas Boilers / Oil Boilers /
# Oil Furnaces (rows 9-11) from placeholder zeros back to the real sales
# figures, which in turn re-computes the dependent fractions on the
# urban-residential / rural-residential / commercial tabs.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the scratch "Texas Notes" sheet entirely.
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete() | Out-Null

# 2. Restore the real sales figures for rows 9-11 on the Data sheet
#    (previously zeroed out / highlighted yellow as placeholders).
$data = $wb.Worksheets.Item("Data")

# Bring back the normal (non-highlighted) number formatting that these
# cells had before they were temporarily zeroed out, by copying the
# formatting from neighboring cells that still use it.
$data.Range("C8").Copy() | Out-Null
$data.Range("C9:C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$data.Range("C14").Copy() | Out-Null
$data.Range("C11").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = $false

$data.Range("C9").Value = 192000
$data.Range("C10").Value = 123000
$data.Range("C11").Value = 56000

# 3. Tidy up sheet selections to reflect where each tab was left.
$about = $wb.Worksheets.Item("About")
$about.Activate()
$about.Range("B17").Select() | Out-Null

$data.Activate()
$data.Range("F30").Select() | Out-Null

$urban = $wb.Worksheets.Item("BFoCSbQL-urban-residential")
$urban.Activate()
$urban.Range("A1").Select() | Out-Null

$rural = $wb.Worksheets.Item("BFoCSbQL-rural-residential")
$rural.Activate()
$rural.Range("A1").Select() | Out-Null

$commercial = $wb.Worksheets.Item("BFoCSbQL-commercial")
$commercial.Activate()
$commercial.Range("A1").Select() | Out-Null

# Leave "About" as the active/selected tab, matching the original workbook.
$about.Activate()

Write-Host "Edit applied."
